$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.215.77"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  +3.40%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.407.73"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  +2.24%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.09%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.86"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.03%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.42"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E7").Value = "  +1.42%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -0.04%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.202"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +10.42%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.02"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +1.85%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000289"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +5.63%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "687.69"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -1.58%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.72"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +3.20%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.952.82"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +1.95%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.106.87"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +3.18%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.403.02"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("E19").Value = "  +1.33%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.40"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +2.17%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.919"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +2.68%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.34"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("E23").Value = "  -0.44%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.81"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  +0.64%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.73"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +1.37%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +2.26%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.96"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  +3.43%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +0.42%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.76"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +11.44%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.17"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +1.17%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "558.02"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("E34").Value = "  +1.63%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.73"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("E36").Value = "  +0.04%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.670.10"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("E38").Value = "  +4.09%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.71"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  +9.43%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.36"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +6.29%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.74"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +4.02%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0430"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +4.74%  "
$ws.Range("E44").Value = "  +1.63%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.69"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +4.53%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -0.30%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.88"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -0.07%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.55"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +2.60%  "
